# Daily attendance processing - 2025-10-30 22:44:29
# Normalize "Recorded By" cell text in column G so that the "System" token
# is listed first, e.g. "user@example.com, System" -> "System, user@example.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Recorded By" value is "dnasr281@gmail.com, System"
$dnasrRows = @(3, 6, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 30, 33, 37, 38, 39, 40, 41, 42, 44, 45, 46, 47, 48, 49, 51, 57, 60, 64, 65, 66, 67, 68, 69, 71, 72, 73, 74, 75, 76, 78, 86, 87, 88, 89, 93, 95, 96, 97, 99, 102, 104, 112, 113, 114, 115, 119, 121, 122, 123, 125, 128, 130, 138, 139, 140, 141, 145, 147, 148, 149, 151, 154, 156)

# Rows whose "Recorded By" value is "admin@admin.com, System"
$adminRows = @(7, 34, 61)

foreach ($r in $dnasrRows) {
    $ws.Range("G" + $r).Value = "System, dnasr281@gmail.com"
}

foreach ($r in $adminRows) {
    $ws.Range("G" + $r).Value = "System, admin@admin.com"
}
